$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 32: add EARNED (C32) of 1.25 ---
$ws.Range("C32").Value = 1.25

# --- Row 33: SL(1-0-0) particular, EARNED 1.25, Absence(H) 1, REMARKS date ---
$ws.Range("B33").Value = "SL(1-0-0)"
$ws.Range("C33").Value = 1.25
$ws.Range("H33").Value = 1
$ws.Range("K31:K32").Copy()
$ws.Range("K33:K34").PasteSpecial(-4122)
$ws.Range("K33").Value = 45224

# --- Row 34: VL(1-0-0) particular, period date, number of days, remarks date ---
$ws.Range("A34").Value = 45231
$ws.Range("B34").Value = "VL(1-0-0)"
$ws.Range("D34").Value = 1
$ws.Range("K34").Value = 45252

# --- Row 35: period date only ---
$ws.Range("A35").Value = 45261

# --- Row 36: year header "2024" (same look as the 2022/2023 rows) ---
$ws.Range("A10").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = "'2024"

# --- Rows 37-48: period dates (monthly) ---
$ws.Range("A37").Value = 45292
$ws.Range("A38").Value = 45323
$ws.Range("A39").Value = 45352
$ws.Range("A40").Value = 45383
$ws.Range("A41").Value = 45413
$ws.Range("A42").Value = 45444
$ws.Range("A43").Value = 45474
$ws.Range("A44").Value = 45505
$ws.Range("A45").Value = 45536
$ws.Range("A46").Value = 45566
$ws.Range("A47").Value = 45597
$ws.Range("A48").Value = 45627

# --- Add a new row to the table (mirrors Table1's "Insert Row Below" for the last row) ---
$ws.Range("A131:K131").Copy()
$ws.Range("A132:K132").PasteSpecial(-4122)
$ws.Range("G132").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("A130:K130").Copy()
$ws.Range("A131:K131").PasteSpecial(-4122)

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K132"))

$ws.Activate()
$ws.Range("K34").Select()
